# First draft of the lab report: diagrams / formatting reshuffled, headers on the
# "actual sheet" tab clarified with the alpha_* error-propagation labels, and the
# "actual sheet" tab becomes the active/selected sheet with a new selection and
# a portrait, A4 page setup (ready for printing as an appendix).

$wb = $excel.ActiveWorkbook

$wsRaw    = $wb.Worksheets.Item(1)   # "not sure what this was"
$wsNicer  = $wb.Worksheets.Item(2)   # "nicer formatting"
$wsActual = $wb.Worksheets.Item(3)   # "actual sheet"

# --- Update header / error-label text on the "actual sheet" tab -----------------
# Order matters only in so far as it controls the order new strings are appended
# to the shared-string table; it mirrors the order the labels appear in the sheet.
$wsActual.Range("J1").Value = "relative error (alpha_R(P))"
$wsActual.Range("C5").Value = "error (alpha_u)"
$wsActual.Range("D5").Value = "relative error (alpha_R(u))"
$wsActual.Range("F5").Value = "error (alpha_(u^-1))"
$wsActual.Range("G5").Value = "relative error (alpha_[R(P)+R(u)])"
$wsActual.Range("J5").Value = "relative error (alpha_[R(P)+R(u)])"
$wsActual.Range("I5").Value = "error (m) (alpha_x)"
$wsActual.Range("I1").Value = "error (alpha_P)"
$wsActual.Range("H1").Value = "pixel size F (m), P"

# --- Make "actual sheet" the active tab with a fresh selection ------------------
[void]$wsActual.Activate()
[void]$wsActual.Range("D6").Select()

# --- Set up the page for this sheet (portrait, A4) ------------------------------
$ps = $wsActual.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
